$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("California")
Write-Host $ws.Name
